$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (order matters for shared-string allocation order)
$ws.Range("BA1").Value = "stevens.p"
$ws.Range("AZ1").Value = "stevens.b"
$ws.Range("BB1").Value = "stevens.r"
$ws.Range("BC1").Value = "mp.b"
$ws.Range("BD1").Value = "mp.p"
$ws.Range("BE1").Value = "mp.r"
$ws.Range("BF1").Value = "quart.b"
$ws.Range("BG1").Value = "quart.p"
$ws.Range("BH1").Value = "quart.r"

# Data values
# Row 2
$ws.Range("AZ2").Value = 0.41586990734776402
$ws.Range("BA2").Value = 0.000000000000000025698035198336401
$ws.Range("BB2").Value = 0.90383766495637996
$ws.Range("BC2").Value = 1.0426166602909499
$ws.Range("BD2").Value = 0.000000029849006417724198
$ws.Range("BE2").Value = 0.84314601571570402
$ws.Range("BF2").Value = 0.13754863532348699
$ws.Range("BG2").Value = 0.00236044354496672
$ws.Range("BH2").Value = 0.42839452582674498
# Row 3
$ws.Range("AZ3").Value = 0.348858294580916
$ws.Range("BA3").Value = 0.0000000000000528512838818042
$ws.Range("BB3").Value = 0.89055601035973098
$ws.Range("BC3").Value = 0.11514084927944
$ws.Range("BD3").Value = 0.61311412481751904
$ws.Range("BE3").Value = 0.015367838413272
$ws.Range("BF3").Value = -0.14013646754939699
$ws.Range("BG3").Value = 0.041492974391440098
$ws.Range("BH3").Value = 0.21119970276444999
# Row 4
$ws.Range("AZ4").Value = -0.046263718581239503
$ws.Range("BA4").Value = 0.40032297457278399
$ws.Range("BB4").Value = 0.032354327225980202
$ws.Range("BC4").Value = -0.0115092535149761
$ws.Range("BD4").Value = 0.92905965520506295
$ws.Range("BE4").Value = 0.00040624451879094097
$ws.Range("BF4").Value = -0.0110955230787968
$ws.Range("BG4").Value = 0.76368122778972802
$ws.Range("BH4").Value = 0.0041942221872483297
# Row 5
$ws.Range("AZ5").Value = 0.0035787540253144302
$ws.Range("BA5").Value = 0.94250203263674504
$ws.Range("BB5").Value = 0.00045182725419334597
$ws.Range("BC5").Value = 0.14741465456706099
$ws.Range("BD5").Value = 0.342454839628338
$ws.Range("BE5").Value = 0.100386206182037
$ws.Range("BF5").Value = -0.0083140289662030006
$ws.Range("BG5").Value = 0.80407182602178096
$ws.Range("BH5").Value = 0.0064483601830224602
# Row 6
$ws.Range("AZ6").Value = 0.19484612780494501
$ws.Range("BA6").Value = 0.0000000000425726378656153
$ws.Range("BB6").Value = 0.71595356615356198
$ws.Range("BC6").Value = 0.17533689866086499
$ws.Range("BD6").Value = 0.054672859670186702
$ws.Range("BE6").Value = 0.15125821648877699
$ws.Range("BF6").Value = -0.020336641438585799
$ws.Range("BG6").Value = 0.31974550962197701
$ws.Range("BH6").Value = 0.032995840210415101
# Row 8
$ws.Range("AZ8").Value = -0.019190487083708501
$ws.Range("BA8").Value = 0.87111480497997895
$ws.Range("BB8").Value = 0.0057973238301720796
$ws.Range("BC8").Value = -0.062734561583077303
$ws.Range("BD8").Value = 0.81931404070529001
$ws.Range("BE8").Value = 0.0061092686968426202
$ws.Range("BF8").Value = 0.077205882352941194
$ws.Range("BG8").Value = 0.28155001721318301
$ws.Range("BH8").Value = 0.082235461945120394
# Row 9
$ws.Range("AZ9").Value = -0.15621454096475801
$ws.Range("BA9").Value = 0.028987107126795
$ws.Range("BB9").Value = 0.64784058836479097
$ws.Range("BC9").Value = -0.15484096571993999
$ws.Range("BD9").Value = 0.27979585944087798
$ws.Range("BE9").Value = 0.096464301932903895
$ws.Range("BF9").Value = -0.0215798825020657
$ws.Range("BG9").Value = 0.50998841948861395
$ws.Range("BH9").Value = 0.023181660699621399
# Row 10
$ws.Range("AZ10").Value = 0.183967229931564
$ws.Range("BA10").Value = 0.00299664279597286
$ws.Range("BB10").Value = 0.79347979439445404
$ws.Range("BC10").Value = 0.00034506556245655299
$ws.Range("BD10").Value = 0.99924644725830503
$ws.Range("BE10").Value = 0.00000025237414669102799
$ws.Range("BF10").Value = 0.056141774891775
$ws.Range("BG10").Value = 0.33323170902379201
$ws.Range("BH10").Value = 0.15577360990063599
# Row 11
$ws.Range("AZ11").Value = 0.049118965362585101
$ws.Range("BA11").Value = 0.41544744799825201
$ws.Range("BB11").Value = 0.0303638340828889
$ws.Range("BC11").Value = 0.30317226890756299
$ws.Range("BD11").Value = 0.216680125273289
$ws.Range("BE11").Value = 0.10679317127647001
$ws.Range("BF11").Value = 0.13804794138770299
$ws.Range("BG11").Value = 0.0064349045259537102
$ws.Range("BH11").Value = 0.30357003453746401
# Row 13
$ws.Range("AZ13").Value = 0.38095009808832198
$ws.Range("BA13").Value = 0.0000000015288364579422501
$ws.Range("BB13").Value = 0.90362274121711805
$ws.Range("BC13").Value = 0.541219256933543
$ws.Range("BD13").Value = 0.070596492103952996
$ws.Range("BE13").Value = 0.24701663211101199
$ws.Range("BF13").Value = 0.120952394249891
$ws.Range("BG13").Value = 0.0162230060488644
$ws.Range("BH13").Value = 0.34765903493539002
# Row 14
$ws.Range("AZ14").Value = 0.53694481569770802
$ws.Range("BA14").Value = 0.0000000000102133222199013
$ws.Range("BB14").Value = 0.91691758045595595
$ws.Range("BC14").Value = 1.6233784013605399
$ws.Range("BD14").Value = 0.0038245045735103302
$ws.Range("BE14").Value = 0.77656202841275301
$ws.Range("BF14").Value = -0.36661904761904801
$ws.Range("BG14").Value = 0.091002277989133198
$ws.Range("BH14").Value = 0.55113015542960597
# Row 15
$ws.Range("AZ15").Value = 0.092904811580424004
$ws.Range("BA15").Value = 0.65544510303558001
$ws.Range("BB15").Value = 0.0751221266857868
$ws.Range("BC15").Value = -0.036590909090909
$ws.Range("BD15").Value = 0.91599820807399102
$ws.Range("BE15").Value = 0.00170510459150112
$ws.Range("BF15").Value = -0.331060606060606
$ws.Range("BG15").Value = 0.0114388979396833
$ws.Range("BH15").Value = 0.62287665692516403
# Row 16
$ws.Range("AZ16").Value = 0.31583780485594698
$ws.Range("BA16").Value = 0.0030752312035444399
$ws.Range("BB16").Value = 0.962836435815698
$ws.Range("BC16").Value = -0.105555555555556
$ws.Range("BD16").Value = 0.72721198827390299
$ws.Range("BE16").Value = 0.018491242824647699
$ws.Range("BF16").Value = 0.058412698412698402
$ws.Range("BG16").Value = 0.364487902659763
$ws.Range("BH16").Value = 0.118428941662146
# Row 17
$ws.Range("AZ17").Value = -0.039664473705920002
$ws.Range("BA17").Value = 0.72447486705119601
$ws.Range("BB17").Value = 0.0079794943386260007
$ws.Range("BC17").Value = -0.17567232385142101
$ws.Range("BD17").Value = 0.45467693864973902
$ws.Range("BE17").Value = 0.047391795199514399
$ws.Range("BF17").Value = -0.0401511000916226
$ws.Range("BG17").Value = 0.46970647693905199
$ws.Range("BH17").Value = 0.035379848534768402
# Row 18
$ws.Range("AZ18").Value = 0.19753394019575199
$ws.Range("BA18").Value = 0.00146400859244754
$ws.Range("BB18").Value = 0.88866651107341799
$ws.Range("BC18").Value = 0.058672253295815803
$ws.Range("BD18").Value = 0.81054514593314297
$ws.Range("BE18").Value = 0.0042448177063020998
$ws.Range("BF18").Value = -0.010726099661540599
$ws.Range("BG18").Value = 0.82330755457079596
$ws.Range("BH18").Value = 0.0026904953189824101
# Row 21
$ws.Range("AZ21").Value = 0.14581039625564801
$ws.Range("BA21").Value = 0.0139586132450954
$ws.Range("BB21").Value = 0.40787072319447598
$ws.Range("BC21").Value = 0.13276075050709901
$ws.Range("BD21").Value = 0.66285647954021498
$ws.Range("BE21").Value = 0.022073426996247699
$ws.Range("BF21").Value = -0.079386294765840207
$ws.Range("BG21").Value = 0.072252958943289702
$ws.Range("BH21").Value = 0.31530394276905899

# Apply scientific-notation number format to very small p-values
$ws.Range("BA2").NumberFormat = "0.00E+00"
$ws.Range("BD2").NumberFormat = "0.00E+00"
$ws.Range("BA3").NumberFormat = "0.00E+00"
$ws.Range("BA6").NumberFormat = "0.00E+00"
$ws.Range("BE10").NumberFormat = "0.00E+00"
$ws.Range("BA13").NumberFormat = "0.00E+00"
$ws.Range("BA14").NumberFormat = "0.00E+00"

# Update selection to match target view state
$ws.Range("AZ8:BH11").Select()
